$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.347.83'
$ws.Range('E2').Value = '  +1.50%  '

$ws.Range('D3').Value = '3.944.92'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.39'
$ws.Range('E5').Value = '  +6.35%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.38'
$ws.Range('E6').Value = '  -1.83%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.615'
$ws.Range('E7').Value = '  -1.52%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.729'
$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('E10').Value = '  -2.59%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000341'
$ws.Range('E11').Value = '  -3.02%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.25'
$ws.Range('E12').Value = '  -2.72%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.44'
$ws.Range('E13').Value = '  -0.27%  '

$ws.Range('D14').Value = '4.577.99'
$ws.Range('E14').Value = '  +0.12%  '

$ws.Range('D15').Value = '3.949.02'
$ws.Range('E15').Value = '  -0.29%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.26'
$ws.Range('E16').Value = '  +6.55%  '

$ws.Range('E17').Value = '  -0.60%  '

$ws.Range('E18').Value = '  +1.69%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.133'
$ws.Range('E19').Value = '  -2.09%  '

$ws.Range('D20').Value = '70.139.46'
$ws.Range('E20').Value = '  +1.11%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '426.03'
$ws.Range('E21').Value = '  -3.07%  '

$ws.Range('E22').Value = '  +1.09%  '

$ws.Range('B23').Value = 'Litecoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '88.36'
$ws.Range('E23').Value = '  -0.68%  '

$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.10'
$ws.Range('E24').Value = '  -3.81%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.02'
$ws.Range('E25').Value = '  +5.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.76'
$ws.Range('E26').Value = '  -3.03%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.63'
$ws.Range('E27').Value = '  -4.71%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.48'
$ws.Range('E28').Value = '  -2.18%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.38'
$ws.Range('E29').Value = '  -0.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '687.84'
$ws.Range('E30').Value = '  -1.91%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.126'
$ws.Range('E31').Value = '  -3.26%  '

$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.87'
$ws.Range('E32').Value = '  -0.57%  '

$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.85'
$ws.Range('E33').Value = '  +12.05%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '67.18'
$ws.Range('E34').Value = '  +7.65%  '

$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0893'
$ws.Range('E35').Value = '  -0.44%  '

$ws.Range('B36').Value = 'TheGraph'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.438'
$ws.Range('E36').Value = '  -5.90%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.75'
$ws.Range('E37').Value = '  -2.98%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.148'
$ws.Range('E38').Value = '  -1.70%  '

$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.38'
$ws.Range('E39').Value = '  +9.71%  '

$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.28%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0479'
$ws.Range('E42').Value = '  -2.16%  '

$ws.Range('E43').Value = '  +4.39%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.81'
$ws.Range('E44').Value = '  -3.64%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.55'
$ws.Range('E45').Value = '  +4.35%  '

$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.19'
$ws.Range('E46').Value = '  +6.31%  '

$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.143'
$ws.Range('E47').Value = '  -0.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000275'
$ws.Range('E48').Value = '  +15.78%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.34'
$ws.Range('E49').Value = '  -1.44%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.97'
$ws.Range('E50').Value = '  +3.32%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0344'
$ws.Range('E51').Value = '  +0.33%  '
